$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - swap B2/C2/D2 values
$ws.Range("B2").Value = 0.7197981723822687
$ws.Range("C2").Value = 0.7197981723822685
$ws.Range("D2").Value = 0.7197981723822685

# Row 3: RandomForestRegressor - updated metrics
$ws.Range("B3").Value = 0.8018740369043749
$ws.Range("C3").Value = 0.8038680273191195
$ws.Range("D3").Value = 0.7632364624274045

# Row 4: model changed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.6914186911079988
$ws.Range("C4").Value = 0.7413077050888451
$ws.Range("D4").Value = 0.7189133595540097

# Row 5: model changed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.6739552204525249
$ws.Range("C5").Value = 0.8274687929408925
$ws.Range("D5").Value = 0.7194445644284359
